$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.432.24"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "2.379.46"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'553.26"
$ws.Range("E5").Value = "  +2.45%  "
$ws.Range("D6").Value = "'139.85"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "2.380.68"
$ws.Range("E9").Value = "  +0.77%  "
$ws.Range("E10").Value = "  +3.16%  "
$ws.Range("D11").Value = "'0.158"
$ws.Range("E11").Value = "  +2.26%  "
$ws.Range("E12").Value = "  +2.57%  "
$ws.Range("D13").Value = "'0.352"
$ws.Range("E13").Value = "  +3.10%  "
$ws.Range("D14").Value = "'25.58"
$ws.Range("E14").Value = "  +3.07%  "
$ws.Range("D15").Value = "'0.0000173"
$ws.Range("E15").Value = "  +6.90%  "
$ws.Range("D16").Value = "2.808.25"
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("D17").Value = "61.415.74"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").Value = "2.378.72"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("E19").Value = "  +3.47%  "
$ws.Range("E20").Value = "  +2.61%  "
$ws.Range("D21").Value = "'321.05"
$ws.Range("E21").Value = "  +1.84%  "
$ws.Range("D22").Value = "'6.70"
$ws.Range("E22").Value = "  +1.80%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  -7.87%  "
$ws.Range("D25").Value = "'64.28"
$ws.Range("E25").Value = "  +1.50%  "
$ws.Range("D26").Value = "'8.85"
$ws.Range("E26").Value = "  +4.53%  "
$ws.Range("E27").Value = "  +0.34%  "
$ws.Range("D28").Value = "2.495.91"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").Value = "'8.17"
$ws.Range("E29").Value = "  +2.67%  "
$ws.Range("D30").Value = "'520.79"
$ws.Range("E30").Value = "  +3.37%  "
$ws.Range("D31").Value = "0.0₃0905"
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("D33").Value = "'0.149"
$ws.Range("E33").Value = "  +2.64%  "
$ws.Range("E34").Value = "  +3.15%  "
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("E37").Value = "  +5.30%  "
$ws.Range("E38").Value = "  +2.93%  "
$ws.Range("E39").Value = "  +6.33%  "
$ws.Range("E40").Value = "  +1.60%  "
$ws.Range("D41").Value = "'18.50"
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("D42").Value = "'146.78"
$ws.Range("E42").Value = "  +6.07%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "'41.38"
$ws.Range("E44").Value = "  +2.95%  "
$ws.Range("D45").Value = "'147.78"
$ws.Range("E45").Value = "  +6.66%  "
$ws.Range("D46").Value = "'2.16"
$ws.Range("E46").Value = "  +2.05%  "
$ws.Range("D47").Value = "'3.61"
$ws.Range("E47").Value = "  +2.73%  "
$ws.Range("E48").Value = "  +2.38%  "
$ws.Range("D49").Value = "'19.77"
$ws.Range("E49").Value = "  +1.62%  "
$ws.Range("D50").Value = "'0.582"
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("E51").Value = "  +1.17%  "
